# Auto-generated edit script: refresh cryptos list values
# (price / volume refresh matching the "Updated cryptos list" GitHub Actions commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.978.16"
$ws.Range("E2").Value = "  +5.49%  "

# Row 3
$ws.Range("D3").Value = "1.917.10"
$ws.Range("E3").Value = "  +4.70%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.22%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.79"
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("E6").Value = "  -0.20%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4748"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.68%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4063"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.69%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.18"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.74%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08181"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.49%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.031"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.34%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.64%  "

# Row 13
$ws.Range("D13").Value = "1.897.79"
$ws.Range("E13").Value = "  +3.92%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.093"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.45%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.395"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.54%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.73%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.19%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001055"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.74%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06632"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.42%  "

# Row 21
$ws.Range("E21").Value = "  -0.26%  "

# Row 22
$ws.Range("D22").Value = "29.008.98"
$ws.Range("E22").Value = "  +5.65%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.589"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.64%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.54%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.270"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.11%  "

# Row 26
$ws.Range("D26").Value = "2.127.48"
$ws.Range("E26").Value = "  +4.14%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.51%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.49%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.191"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.99%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.543"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.71%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.19%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.013"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.59%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09583"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.00%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.441"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.02%  "

# Row 35
$ws.Range("E35").Value = "  +1.73%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.438"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.56%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06210"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.69%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.708"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.67%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02284"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.42%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.202"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.45%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6051"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.66%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.06%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1902"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.00%  "

# Row 44
$ws.Range("E44").Value = "  -0.18%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.281"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.30%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5643"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.36%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.83%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.984"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.01%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07248"
$ws.Range("D49").Style = "Normal"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.156"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +18.27%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.19%  "
